$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 81891
$ws.Range("B2").Value = "Sra. Lavínia Almeida"
$ws.Range("C2").Value = "Financeiro"
$ws.Range("D2").Value = "Consulta medica"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 45103
$ws.Range("G2").Value = 5916.7

# Row 3
$ws.Range("A3").Value = 98527
$ws.Range("B3").Value = "João Vitor Montenegro"
$ws.Range("C3").Value = "Vendas"
$ws.Range("D3").Value = "Problemas pessoais"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 45100
$ws.Range("G3").Value = 6702.45

# Row 4
$ws.Range("A4").Value = 3544
$ws.Range("B4").Value = "Yan da Rocha"
$ws.Range("C4").Value = "Atendimento ao Cliente"
$ws.Range("D4").Value = "Problemas pessoais"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 45080
$ws.Range("G4").Value = 5564.77

# Row 5
$ws.Range("A5").Value = 40353
$ws.Range("B5").Value = "Luiz Henrique Macedo"
$ws.Range("C5").Value = "P&D"
$ws.Range("D5").Value = "Consulta medica"
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 45104
$ws.Range("G5").Value = 7391.76

# Row 6
$ws.Range("A6").Value = 82836
$ws.Range("B6").Value = "Dr. Carlos Eduardo Andrade"
$ws.Range("C6").Value = "Marketing"
$ws.Range("D6").Value = "Problemas pessoais"
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 45081
$ws.Range("G6").Value = 9871.51

# Row 7
$ws.Range("A7").Value = 4751
$ws.Range("B7").Value = "Ágatha Costa"
$ws.Range("C7").Value = "Operacoes"
$ws.Range("D7").Value = "Problemas pessoais"
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 45088
$ws.Range("G7").Value = 3541.27

# Row 8
$ws.Range("A8").Value = 56498
$ws.Range("B8").Value = "João Felipe Costa"
$ws.Range("C8").Value = "Financeiro"
$ws.Range("D8").Value = "Outros"
$ws.Range("E8").Value = 6
$ws.Range("F8").Value = 45106
$ws.Range("G8").Value = 8030.39

# Row 9
$ws.Range("A9").Value = 59547
$ws.Range("B9").Value = "Heloisa Peixoto"
$ws.Range("C9").Value = "Engenharia"
$ws.Range("D9").Value = "Problemas pessoais"
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = 45091
$ws.Range("G9").Value = 3918.06

# Row 10
$ws.Range("A10").Value = 85339
$ws.Range("B10").Value = "Dr. Brayan Cirino"
$ws.Range("C10").Value = "TI"
$ws.Range("D10").Value = "Doenca"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 45103
$ws.Range("G10").Value = 7104.82

# Row 11
$ws.Range("A11").Value = 36915
$ws.Range("B11").Value = "Alana da Mota"
$ws.Range("C11").Value = "Juridico"
$ws.Range("D11").Value = "Viagem de negocios"
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 45096
$ws.Range("G11").Value = 5779.85
